$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data feed rows appended to the Source/Parameter/URL/Notes table:
#   Row 9  - EIA / E85 / afdc.energy.gov price source / download-button note
#   Row 10 - Canpotex / Freight Cost / Ben's historical reference file / sum note
# Cells are written in the same left-to-right, row-by-row order the rows were
# authored in so new shared-string entries land in the same order.

$ws.Range("A9").Value = "EIA"
$ws.Range("B9").Value = "E85"
$ws.Range("C9").Value = "https://afdc.energy.gov/fuels/prices.html"
$ws.Range("D9").Value = "E85 >> Download button"

$ws.Range("A10").Value = "Canpotex"
$ws.Range("D10").Value = "N to Q and AK >> Sum"
$ws.Range("B10").Value = "Freight Cost"
$ws.Range("C10").Value = 'From Ben - Historical File "Reference: Monthly Netback Testing.xlsx"'

# Move the active selection to the next empty row, matching the workbook's
# post-edit UI state.
$ws.Range("A11").Select()
